$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date stamp in A1 (was 45310 -> 45311, i.e. +1 day)
$ws.Range("A1").Value = 45311

# Update prices in the price list
$ws.Range("D30").Value = 134
$ws.Range("D31").Value = 144.78
